$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich-text shared strings): bump the report's
# volume/number and week-covering dates by one week.
# ---------------------------------------------------------------------

# "Volume 29   Number  45" -> "...46"  (characters 21-22 are "45")
$ws.Range("A8").Characters(21, 2).Text = "46"

# "Report Covering the Week  11/7/2022  Through  11/13/2022"
#   11/7/2022  (chars 27-35, 9 chars) -> 11/14/2022
#   11/13/2022 (chars 48-57 after the above edit, 10 chars) -> 11/20/2022
$ws.Range("C9").Characters(27, 9).Text = "11/14/2022"
$ws.Range("C9").Characters(48, 10).Text = "11/20/2022"

# ---------------------------------------------------------------------
# Row 15 - only one figure changed
# ---------------------------------------------------------------------
$ws.Range("L15").Value = -16.666666666666

# ---------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 41.666666666666
$ws.Range("I16").Value = 216
$ws.Range("J16").Value = 158
$ws.Range("K16").Value = 36.708860759493
$ws.Range("L16").Value = 3.349282296650
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = -82.119205298013

# ---------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -46.153846153846
$ws.Range("I17").Value = 149
$ws.Range("J17").Value = 144
$ws.Range("K17").Value = 3.472222222222
$ws.Range("L17").Value = 40.566037735849
$ws.Range("M17").Value = 73.255813953488
$ws.Range("N17").Value = -43.984962406015

# ---------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 225
$ws.Range("J18").Value = 207
$ws.Range("K18").Value = 8.695652173913
$ws.Range("M18").Value = 3.211009174311
$ws.Range("N18").Value = -91.499811106913

# ---------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 34
$ws.Range("E19").Value = -29.411764705882
$ws.Range("F19").Value = 130
$ws.Range("G19").Value = 129
$ws.Range("H19").Value = 0.775193798449
$ws.Range("I19").Value = 1551
$ws.Range("J19").Value = 1046
$ws.Range("K19").Value = 48.279158699808
$ws.Range("L19").Value = 30.008382229673
$ws.Range("M19").Value = 38.482142857142
$ws.Range("N19").Value = -53.409432261940

# ---------------------------------------------------------------------
# Row 20 - C/F/G/H/I/K/L/M/N numeric edits; D20 & E20 switch from
# numbers to the "0" / "***.*" placeholder text used elsewhere in the
# sheet, so copy format+value from cells that already hold that
# placeholder text before overwriting the numeric cells below.
# ---------------------------------------------------------------------
$ws.Range("C15").Copy($ws.Range("D20"))
$ws.Range("E15").Copy($ws.Range("E20"))

$ws.Range("C20").Value = 6
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 36.363636363636
$ws.Range("I20").Value = 162
$ws.Range("K20").Value = 5.882352941176
$ws.Range("L20").Value = 45.945945945945
$ws.Range("M20").Value = 95.180722891566
$ws.Range("N20").Value = -94.747081712062

# ---------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 185
$ws.Range("G21").Value = 188
$ws.Range("H21").Value = -1.595744680851
$ws.Range("I21").Value = 2315
$ws.Range("J21").Value = 1722
$ws.Range("K21").Value = 34.436701509872
$ws.Range("L21").Value = 16.214859437751
$ws.Range("M21").Value = 41.072516758074
$ws.Range("N21").Value = -78.094246782740

# ---------------------------------------------------------------------
# Row 22 - C22 switches from the "0" placeholder text to a real number.
# ---------------------------------------------------------------------
$ws.Range("D22").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2

$ws.Range("E22").Value = 100
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 43
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 186.666666666667
$ws.Range("L22").Value = 43.333333333333
$ws.Range("M22").Value = 79.166666666666

# ---------------------------------------------------------------------
# Row 23 - C23 and G23 switch to the "0" placeholder text, H23 switches
# to the "***.*" placeholder text.
# ---------------------------------------------------------------------
$ws.Range("D23").Copy($ws.Range("C23"))
$ws.Range("D23").Copy($ws.Range("G23"))
$ws.Range("E23").Copy($ws.Range("H23"))

$ws.Range("I23").Value = 33
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = 6.451612903225
$ws.Range("M23").Value = 37.5

# ---------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 104
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 116.666666666667
$ws.Range("F24").Value = 372
$ws.Range("G24").Value = 180
$ws.Range("H24").Value = 106.666666666667
$ws.Range("I24").Value = 3603
$ws.Range("J24").Value = 2182
$ws.Range("K24").Value = 65.123739688359
$ws.Range("L24").Value = 52.734209410767
$ws.Range("M24").Value = 138.609271523179

# ---------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 20.833333333333
$ws.Range("I25").Value = 340
$ws.Range("J25").Value = 305
$ws.Range("K25").Value = 11.475409836065
$ws.Range("L25").Value = 37.651821862348
$ws.Range("M25").Value = 8.280254777070

# ---------------------------------------------------------------------
# Row 26 - D26 and E26 switch from placeholder text to real numbers.
# ---------------------------------------------------------------------
$ws.Range("G26").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1

$ws.Range("H26").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100

$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -30
$ws.Range("L26").Value = -17.647058823529

# ---------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("I27").Value = 74
$ws.Range("J27").Value = 85
$ws.Range("K27").Value = -12.941176470588
$ws.Range("L27").Value = 23.333333333333
